# typekeys.xlsx update
#
# [json] add storeKeys(json,jsonpath,var) - extract immediate keys of a
# resolved JSON fragment based on jsonpath. The "#system" hidden sheet
# keeps, per-category, an alphabetically sorted list of function
# signatures. Adding "storeKeys" to the "json" category means:
#   - a new row is inserted (alphabetically, right before storeValue)
#     in the json column (M) on the "#system" sheet, and
#   - the "json" defined name's range grows by one row.
#
# Also: the "text" category (a single-row category sitting in column Y)
# is removed, which shifts every category column from Z onward one
# column to the left (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD), and
# shifts the "text"/"web"/.../"xml" category-name rows in column A
# (the "target" defined name list) up by one row as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) "json" category (column M): insert storeKeys(json,jsonpath,var)
#    right before storeValue(json,jsonpath,var), i.e. at M16, pushing
#    storeValue/storeValues down by one row (M16:M17 -> M17:M18).
# ---------------------------------------------------------------------
$jsonTail = $ws.Range("M16:M17").Value2
$ws.Range("M17:M18").Value2 = $jsonTail
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------
# 2) "target" category names (column A): the "text" category (A25) goes
#    away, so web/webalert/webcookie/ws/ws.async/xml (A26:A31) shift up
#    to A25:A30, and the old last row (A31) is cleared.
# ---------------------------------------------------------------------
$targetTail = $ws.Range("A26:A31").Value2
$ws.Range("A25:A30").Value2 = $targetTail
$ws.Range("A31").Clear()

# ---------------------------------------------------------------------
# 3) category data columns: web/webalert/webcookie/ws/ws.async/xml each
#    shift one column to the left (their old "text" neighbor in Y is
#    gone). Read each source column before overwriting, then shift
#    left-to-right so every read happens before that column is
#    overwritten.
# ---------------------------------------------------------------------
$web       = $ws.Range("Z1:Z129").Value2
$webalert  = $ws.Range("AA1:AA129").Value2
$webcookie = $ws.Range("AB1:AB129").Value2
$wsCol     = $ws.Range("AC1:AC129").Value2
$wsAsync   = $ws.Range("AD1:AD129").Value2
$xml       = $ws.Range("AE1:AE129").Value2

$ws.Range("Y1:Y129").Value2  = $web
$ws.Range("Z1:Z129").Value2  = $webalert
$ws.Range("AA1:AA129").Value2 = $webcookie
$ws.Range("AB1:AB129").Value2 = $wsCol
$ws.Range("AC1:AC129").Value2 = $wsAsync
$ws.Range("AD1:AD129").Value2 = $xml
$ws.Range("AE1:AE129").Clear()

# ---------------------------------------------------------------------
# 4) defined names: keep every range in sync with the shifted data.
# ---------------------------------------------------------------------
$wb.Names.Item("json").RefersTo       = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo     = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo        = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AD`$2:`$AD`$27"
